$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column G (Return_with_prediction) for rows 2-57
$ws.Range("G2").Value = 0.07175570646056607
$ws.Range("G3").Value = 0.08579977571796248
$ws.Range("G4").Value = 0.01285023862399317
$ws.Range("G5").Value = 0.03667978134259088
$ws.Range("G6").Value = -0.1391201148019596
$ws.Range("G7").Value = -0.100551325032716
$ws.Range("G8").Value = -0.183368182132807
$ws.Range("G9").Value = -0.324975731188394
$ws.Range("G10").Value = 0.001629175921530999
$ws.Range("G11").Value = 0.04036417612975503
$ws.Range("G12").Value = 0.1904457253719538
$ws.Range("G13").Value = 0.2188362450285577
$ws.Range("G14").Value = -0.05660425079313351
$ws.Range("G15").Value = -0.04037792320375548
$ws.Range("G16").Value = 0.1735460257360774
$ws.Range("G17").Value = 0.1967038881324487
$ws.Range("G18").Value = 0.06167126440519453
$ws.Range("G19").Value = 0.07253740722019078
$ws.Range("G20").Value = 0.01174711753588197
$ws.Range("G21").Value = -0.01576337121632567
$ws.Range("G22").Value = 0.0551671531278437
$ws.Range("G23").Value = 0.05919744231952283
$ws.Range("G24").Value = 0.03055488235836543
$ws.Range("G25").Value = 0.01407829811901563
$ws.Range("G26").Value = 0.1214295611502891
$ws.Range("G27").Value = 0.1335262539901777
$ws.Range("G28").Value = 0.1069109252595103
$ws.Range("G29").Value = 0.1338838393465525
$ws.Range("G30").Value = 0.05603510954370651
$ws.Range("G31").Value = 0.07033973392471024
$ws.Range("G32").Value = 0.079023832728689
$ws.Range("G33").Value = 0.079023832728689
$ws.Range("G34").Value = 0.01047382111804239
$ws.Range("G35").Value = 0.02997738422266027
$ws.Range("G36").Value = -0.008973995070044096
$ws.Range("G37").Value = 0.004457551931185209
$ws.Range("G38").Value = 0.0622025812485167
$ws.Range("G39").Value = 0.03849974700149494
$ws.Range("G40").Value = 0.06282465743436066
$ws.Range("G41").Value = 0.07710818430837971
$ws.Range("G42").Value = 0.0535117737405689
$ws.Range("G43").Value = 0.06941806142794307
$ws.Range("G44").Value = 0.1004535807871566
$ws.Range("G45").Value = 0.1091476905107295
$ws.Range("G46").Value = 0.01849816019712222
$ws.Range("G47").Value = -0.004335443576684862
$ws.Range("G48").Value = -0.009835836566473627
$ws.Range("G49").Value = 0.01222607316051134
$ws.Range("G50").Value = 0.1402935395337438
$ws.Range("G51").Value = 0.1550049061531065
$ws.Range("G52").Value = 0.07965484767475919
$ws.Range("G53").Value = 0.06837753631612511
$ws.Range("G54").Value = -0.1069902304093575
$ws.Range("G55").Value = -0.09731671039542743
$ws.Range("G56").Value = 0.1689111944398531
$ws.Range("G57").Value = 0.1648522972030135

# Update column H (return_pct_change) for rows 2-57
$ws.Range("H2").Value = 115.5151085434947
$ws.Range("H3").Value = 70.01939830953505
$ws.Range("H4").Value = -73.60036646361013
$ws.Range("H5").Value = -48.08574424034298
$ws.Range("H6").Value = -18.27359133040317
$ws.Range("H7").Value = 19.58047331513407
$ws.Range("H8").Value = 7.937866628365454
$ws.Range("H9").Value = -7.580468739299771
$ws.Range("H10").Value = 28.00791023623183
$ws.Range("H11").Value = 296.2915011472496
$ws.Range("H12").Value = -10.12949973140133
$ws.Range("H13").Value = -6.876842088008631
$ws.Range("H14").Value = 37.83186891356716
$ws.Range("H15").Value = 43.10462934326742
$ws.Range("H16").Value = -9.353328016551783
$ws.Range("H17").Value = 13.1329811389915
$ws.Range("H18").Value = 13.72197653273247
$ws.Range("H19").Value = -15.69299935289367
$ws.Range("H20").Value = -7.713074770338974
$ws.Range("H21").Value = 70.74018478423035
$ws.Range("H22").Value = -15.48209393465643
$ws.Range("H23").Value = 2.642558945501176
$ws.Range("H24").Value = -5.677661473116837
$ws.Range("H25").Value = -52.16690762190629
$ws.Range("H26").Value = 7.183374741995325
$ws.Range("H27").Value = 48.0537999798629
$ws.Range("H28").Value = -8.986507781131444
$ws.Range("H29").Value = 11.91712209339007
$ws.Range("H30").Value = -16.64956340466355
$ws.Range("H31").Value = 2.509662683723797
$ws.Range("H32").Value = 80.97106811356278
$ws.Range("H33").Value = 45.42904817693142
$ws.Range("H34").Value = 154.8407266748787
$ws.Range("H35").Value = 114.8505433452264
$ws.Range("H36").Value = -158.0506047770673
$ws.Range("H37").Value = -64.40698765535876
$ws.Range("H38").Value = -13.29780618508528
$ws.Range("H39").Value = -10.59566517347874
$ws.Range("H40").Value = 40.45380144796464
$ws.Range("H41").Value = 523.8868346606756
$ws.Range("H42").Value = 2.35810917625428
$ws.Range("H43").Value = 39.11467136015339
$ws.Range("H44").Value = -23.75489904603738
$ws.Range("H45").Value = -39.1686416415331
$ws.Range("H46").Value = 142.1084150170526
$ws.Range("H47").Value = -65.50231102169892
$ws.Range("H48").Value = -167.8716410674847
$ws.Range("H49").Value = 319.9401419210821
$ws.Range("H50").Value = -1.845591515295671
$ws.Range("H51").Value = 18.35578373806961
$ws.Range("H52").Value = 28.5736410636679
$ws.Range("H53").Value = 11.76155016647274
$ws.Range("H54").Value = -19.80057830663565
$ws.Range("H55").Value = 6.233848474900645
$ws.Range("H56").Value = 8.888117314086101
$ws.Range("H57").Value = 18.18546805484503

# Update column I (mean_return_pct_change) for row 2
$ws.Range("I2").Value = 26.492478762591
